$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct power band calculation for row 23 (B23: 6 -> 5)
$ws.Range("B23").Value = 5

# Add two new rows of data (24 and 25)
$ws.Range("A24").Value = 1061
$ws.Range("B24").Value = 5
$ws.Range("C24").Value = 2
$ws.Range("D24").Value = 0.48321759259259256
$ws.Range("D24").NumberFormat = $ws.Range("D23").NumberFormat

$ws.Range("A25").Value = 1064
$ws.Range("B25").Value = 5
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 0.48321759259259256
$ws.Range("D25").NumberFormat = $ws.Range("D23").NumberFormat

# Update selection to match target state
$ws.Range("D26").Select()
